# "Added a quick and dirty Change over time"
#
# Fixes a couple of existing values in the storm-interval table and
# appends five new rows (116-120) describing later storms, also filling
# in the "next storm start/end" columns (E/F) for rows 113-115 which
# previously pointed nowhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "YYYY-MM-DD HH:MM:SS"

function Set-DateCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $value
    $c.NumberFormat = $dateFormat
}

# Re-assert the date/time display format on every existing start/end cell
# (columns B, C, E, F) so the whole table keeps a single, consistent date
# style after the edit.
$ws.Range("B2:C115").NumberFormat = $dateFormat
$ws.Range("E2:F112").NumberFormat = $dateFormat

function Set-IndexCell($row, $value) {
    $c = $ws.Cells.Item($row, 1)
    $c.Value = $value
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108  # xlCenter
    $c.VerticalAlignment = -4160    # xlTop
    $c.Borders.LineStyle = 1
}

# --- Fix existing values -------------------------------------------------

# Row 57: next-storm-start (E57) corrected.
Set-DateCell 57 5 41684.625

# Row 60: start (B60) corrected to match, with duration (D60) and the
# resulting duration-in-ticks (G60) recomputed accordingly.
Set-DateCell 60 2 41684.625
$ws.Cells.Item(60, 4).Value = 5.25
$ws.Cells.Item(60, 7).Value = 18900000000000

# --- Fill in "next storm" start/end for rows 113-115 ---------------------

Set-DateCell 113 5 42005.354166666664
Set-DateCell 113 6 42005.895833333336

Set-DateCell 114 5 42006.5625
Set-DateCell 114 6 42006.854166666664

Set-DateCell 115 5 42008.8125
Set-DateCell 115 6 42008.947916666664

# --- Append new rows 116-120 ---------------------------------------------

$newRows = @(
    @{ Row = 116; A = 122; B = 42005.354166666664; C = 42005.895833333336; D = 13.0;  E = 42009.125;            F = 42009.364583333336; G = 46800000000000 },
    @{ Row = 117; A = 123; B = 42006.5625;          C = 42006.854166666664; D = 7.0;   E = 42013.489583333336;  F = 42014.375;           G = 25200000000000 },
    @{ Row = 118; A = 124; B = 42008.8125;          C = 42008.947916666664; D = 3.25;  E = $null;               F = $null;               G = 11700000000000 },
    @{ Row = 119; A = 125; B = 42009.125;           C = 42009.364583333336; D = 5.75;  E = $null;               F = $null;               G = 20700000000000 },
    @{ Row = 120; A = 126; B = 42013.489583333336;  C = 42014.375;          D = 21.25; E = $null;               F = $null;               G = 76500000000000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    Set-IndexCell $row $r.A

    Set-DateCell $row 2 $r.B
    Set-DateCell $row 3 $r.C

    $ws.Cells.Item($row, 4).Value = $r.D

    if ($null -ne $r.E) {
        Set-DateCell $row 5 $r.E
    }
    if ($null -ne $r.F) {
        Set-DateCell $row 6 $r.F
    }

    $ws.Cells.Item($row, 7).Value = $r.G
}
